# Auto-generated Excel COM-interop script to apply the data refresh diff
# to Sheets/Behemoth_Profits.xlsx (workbook with sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 373.9
$ws.Range("I9").Value = 346.125
$ws.Range("J9").Value = 485
$ws.Range("K9").Value = 346.125
$ws.Range("L9").Value = 485
$ws.Range("M9").Value = -177.125
$ws.Range("N9").Value = -823

$ws.Range("H44").Value = 17833
$ws.Range("J44").Value = 17833
$ws.Range("L44").Value = 17833
$ws.Range("N44").Value = -18757

$ws.Range("H64").Value = 4905.1577
$ws.Range("I64").Value = 4433
$ws.Range("K64").Value = 4433
$ws.Range("M64").Value = -4185

$ws.Range("H67").Value = 4905.1577
$ws.Range("I67").Value = 4433
$ws.Range("K67").Value = 4433
$ws.Range("M67").Value = -3575

$ws.Range("H74").Value = 2999
$ws.Range("I74").Value = 2999
$ws.Range("K74").Value = 2999
$ws.Range("M74").Value = -2063

$ws.Range("H77").Value = 2999
$ws.Range("I77").Value = 2999
$ws.Range("K77").Value = 14995
$ws.Range("M77").Value = -10315

$ws.Range("H112").Value = 1872.3334
$ws.Range("I112").Value = 1076.6666
$ws.Range("J112").Value = 2031.4667
$ws.Range("K112").Value = 3229.9998
$ws.Range("L112").Value = 6094.4001
$ws.Range("M112").Value = -2121.9998
$ws.Range("N112").Value = -8310.400099999999

$ws.Range("H125").Value = 1906.25
$ws.Range("I125").Value = 1190.2727
$ws.Range("K125").Value = 10712.4543
$ws.Range("M125").Value = -8252.454299999999

$ws.Range("H138").Value = 2344.4614
$ws.Range("I138").Value = 1209.0588
$ws.Range("J138").Value = 2895.9429
$ws.Range("K138").Value = 3627.1764
$ws.Range("L138").Value = 8687.8287
$ws.Range("M138").Value = 1512.8236
$ws.Range("N138").Value = -18967.8287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 41333.332
$ws.Range("J7").Value = 52000
$ws.Range("L7").Value = 52000
$ws.Range("N7").Value = -52228

$ws.Range("H32").Value = 9641949
$ws.Range("I32").Value = 11390394
$ws.Range("K32").Value = 11390394
$ws.Range("M32").Value = -11390107

$ws.Range("H45").Value = 20002294
$ws.Range("I45").Value = 31251538
$ws.Range("J45").Value = 3634.7778
$ws.Range("K45").Value = 31251538
$ws.Range("L45").Value = 3634.7778
$ws.Range("M45").Value = -31251161
$ws.Range("N45").Value = -4388.7778

$ws.Range("H63").Value = 5623.4614
$ws.Range("I63").Value = 2183.3333
$ws.Range("J63").Value = 8572.143
$ws.Range("K63").Value = 2183.3333
$ws.Range("L63").Value = 8572.143
$ws.Range("M63").Value = -1497.3333
$ws.Range("N63").Value = -9944.143

$ws.Range("H66").Value = 5623.4614
$ws.Range("I66").Value = 2183.3333
$ws.Range("J66").Value = 8572.143
$ws.Range("K66").Value = 10916.6665
$ws.Range("L66").Value = 42860.715
$ws.Range("M66").Value = -7484.666499999999
$ws.Range("N66").Value = -49724.715

$ws.Range("H68").Value = 39990
$ws.Range("I68").Value = 39990
$ws.Range("K68").Value = 39990
$ws.Range("M68").Value = -39179

$ws.Range("H71").Value = 39990
$ws.Range("I71").Value = 39990
$ws.Range("K71").Value = 119970
$ws.Range("M71").Value = -115914

$ws.Range("H110").Value = 1666.5518
$ws.Range("I110").Value = 1584.2916
$ws.Range("K110").Value = 1584.2916
$ws.Range("M110").Value = 460.7084

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 1200
$ws.Range("I31").Value = 1200
$ws.Range("K31").Value = 1200
$ws.Range("M31").Value = -948

$ws.Range("H62").Value = 128000
$ws.Range("J62").Value = 128000
$ws.Range("L62").Value = 128000
$ws.Range("N62").Value = -129372

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = ""
$ws.Range("N63").Value = 0

$ws.Range("H65").Value = 128000
$ws.Range("J65").Value = 128000
$ws.Range("L65").Value = 384000
$ws.Range("N65").Value = -390864

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = ""
$ws.Range("N66").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1705.75
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1705.75
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = 1705.75
$ws.Range("N16").Value = -2279.75

$ws.Range("H26").Value = 13380.77
$ws.Range("J26").Value = 13880
$ws.Range("L26").Value = 13880
$ws.Range("N26").Value = -14454

$ws.Range("H58").Value = 2094.9395
$ws.Range("I58").Value = 1231.5385
$ws.Range("J58").Value = 5301.857
$ws.Range("K58").Value = 1231.5385
$ws.Range("L58").Value = 5301.857
$ws.Range("M58").Value = -1028.5385
$ws.Range("N58").Value = -5707.857

$ws.Range("H74").Value = 37708.43
$ws.Range("J74").Value = 37708.43
$ws.Range("L74").Value = 37708.43
$ws.Range("N74").Value = -39456.43

$ws.Range("H77").Value = 37708.43
$ws.Range("J77").Value = 37708.43
$ws.Range("L77").Value = 113125.29
$ws.Range("N77").Value = -121861.29

$ws.Range("H87").Value = 116979
$ws.Range("J87").Value = 116979
$ws.Range("L87").Value = 116979
$ws.Range("N87").Value = -119351

$ws.Range("H90").Value = 116979
$ws.Range("J90").Value = 116979
$ws.Range("L90").Value = 350937
$ws.Range("N90").Value = -362793

$ws.Range("H105").Value = 999.86957
$ws.Range("I105").Value = 993.4737
$ws.Range("J105").Value = 1030.25
$ws.Range("K105").Value = 993.4737
$ws.Range("L105").Value = 1030.25
$ws.Range("M105").Value = 753.5263
$ws.Range("N105").Value = -4524.25

$ws.Range("H113").Value = 1705.75
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1705.75
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = ""
$ws.Range("M113").Value = 1705.75
$ws.Range("N113").Value = -6045.75

$ws.Range("H136").Value = 2094.9395
$ws.Range("I136").Value = 1231.5385
$ws.Range("J136").Value = 5301.857
$ws.Range("K136").Value = 3694.6155
$ws.Range("L136").Value = 15905.571
$ws.Range("M136").Value = -1144.6155
$ws.Range("N136").Value = -21005.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 168.14285
$ws.Range("J26").Value = 124.25
$ws.Range("L26").Value = 372.75
$ws.Range("N26").Value = -948.75

$ws.Range("H29").Value = 103.1
$ws.Range("J29").Value = 133.66667
$ws.Range("L29").Value = 401.00001
$ws.Range("N29").Value = -955.00001

$ws.Range("H31").Value = 574.75
$ws.Range("I31").Value = 550
$ws.Range("J31").Value = 599.5
$ws.Range("K31").Value = 1650
$ws.Range("L31").Value = 1798.5
$ws.Range("M31").Value = -1362
$ws.Range("N31").Value = -2374.5

$ws.Range("H51").Value = 23002
$ws.Range("J51").Value = 25670
$ws.Range("L51").Value = 77010
$ws.Range("N51").Value = -77930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 24999
$ws.Range("J28").Value = 24999
$ws.Range("L28").Value = 24999
$ws.Range("N28").Value = -25383

$ws.Range("H80").Value = 2999.6924
$ws.Range("I80").Value = 2719
$ws.Range("K80").Value = 2719
$ws.Range("M80").Value = -1721

$ws.Range("H83").Value = 2999.6924
$ws.Range("I83").Value = 2719
$ws.Range("K83").Value = 13595
$ws.Range("M83").Value = -8603

$ws.Range("H100").Value = 57784.668
$ws.Range("J100").Value = 57784.668
$ws.Range("L100").Value = 57784.668
$ws.Range("N100").Value = -59948.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 99270.17999999999
$ws.Range("I7").Value = 6613.25
$ws.Range("K7").Value = 6613.25
$ws.Range("M7").Value = -6501.25

$ws.Range("H100").Value = 3493.5
$ws.Range("I100").Value = 2384.6
$ws.Range("K100").Value = 2384.6
$ws.Range("M100").Value = -1843.6

$ws.Range("H126").Value = 99270.17999999999
$ws.Range("I126").Value = 6613.25
$ws.Range("K126").Value = 19839.75
$ws.Range("M126").Value = -17369.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 5000
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 5000
$ws.Range("N34").Value = -5406

$ws.Range("H37").Value = 15547.5
$ws.Range("J37").Value = 15547.5
$ws.Range("L37").Value = 15547.5
$ws.Range("N37").Value = -15953.5

$ws.Range("H107").Value = 1432.1052
$ws.Range("I107").Value = 1345.9286
$ws.Range("J107").Value = 1673.4
$ws.Range("K107").Value = 4037.7858
$ws.Range("L107").Value = 5020.200000000001
$ws.Range("M107").Value = -2117.7858
$ws.Range("N107").Value = -8860.200000000001

Write-Host "Applied Behemoth_Profits data refresh."